$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "DB1": update the selection (no data change)
# ---------------------------------------------------------------
$wsDB1 = $wb.Worksheets.Item("DB1")
$wsDB1.Activate()
$wsDB1.Range("A2:E2").Select()

# ---------------------------------------------------------------
# Sheet "DB3": insert a new STRUCT row after the header, and add a
# trailing END_STRUCT row - mirroring the pattern already used on
# sheets DB2 and DB1 (read/write support added for dword, byte,
# word, dint, time, date, time of day and char).
# ---------------------------------------------------------------
$wsDB3 = $wb.Worksheets.Item("DB3")
$wsDB3.Activate()

# Insert a new row 2 (pushes existing rows 2.. down by one)
$wsDB3.Rows.Item(2).Insert()

# Row 2 becomes the STRUCT opener, matching DB2/DB1 layout
$wsDB3.Range("A2").NumberFormat = "@"
$wsDB3.Range("A2").Value2 = "0.0"
$wsDB3.Range("C2").Value2 = "STRUCT"

# Append the new closing END_STRUCT row (row 13)
$wsDB3.Range("A13").NumberFormat = "@"
$wsDB3.Range("A13").Value2 = "28.0"
$wsDB3.Range("C13").Value2 = "END_STRUCT"

$wsDB3.Range("F11").Select()

Write-Host "done"
